# Update the "single-step conv-lstm" connector so it no longer terminates
# on the small oval's connection site and instead runs down to the
# markdown box below it (per the commit "added markdown in
# single-step_conv-lstm, update").

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$connector = $s.Shapes.Item("Straight Arrow Connector 85")

# Detach the arrow's end point from shape 77 (it used to snap to that
# oval's connection site idx=2); leave the start point attached to shape 39.
$connector.ConnectorFormat.EndDisconnect()

# Stretch the connector's end point down (its start position/anchor at
# shape 39 is unchanged) so it now reaches the new markdown content
# added below, instead of stopping at the oval's connection site.
$connector.Height = 2297151 / 12700
